# Update "想去人数" (want-to-go count) figures in column F across all
# four sheets of the workbook, matching the refreshed scrape snapshot.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 247
$ws1.Range("F7").Value = 12921
$ws1.Range("F10").Value = 232
$ws1.Range("F11").Value = 2985
$ws1.Range("F12").Value = 80
$ws1.Range("F13").Value = 6337
$ws1.Range("F16").Value = 3358
$ws1.Range("F17").Value = 25
$ws1.Range("F20").Value = 34
$ws1.Range("F23").Value = 23
$ws1.Range("F24").Value = 3573
$ws1.Range("F25").Value = 81
$ws1.Range("F27").Value = 2723
$ws1.Range("F28").Value = 2723
$ws1.Range("F29").Value = 397
$ws1.Range("F30").Value = 1864
$ws1.Range("F31").Value = 99
$ws1.Range("F33").Value = 6537
$ws1.Range("F36").Value = 567
$ws1.Range("F37").Value = 1959
$ws1.Range("F39").Value = 88
$ws1.Range("F40").Value = 1019
$ws1.Range("F43").Value = 216
$ws1.Range("F46").Value = 123
$ws1.Range("F47").Value = 1186
$ws1.Range("F48").Value = 1742

# Sheet "演出" (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F15").Value = 95

# Sheet "本地生活" (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 580

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 580
$ws4.Range("F8").Value = 247
$ws4.Range("F10").Value = 12921
$ws4.Range("F13").Value = 232
$ws4.Range("F14").Value = 2985
$ws4.Range("F15").Value = 6337
$ws4.Range("F17").Value = 3358
$ws4.Range("F18").Value = 25
$ws4.Range("F21").Value = 34
$ws4.Range("F25").Value = 23
$ws4.Range("F26").Value = 3573
$ws4.Range("F28").Value = 2723
$ws4.Range("F29").Value = 397
$ws4.Range("F30").Value = 1864
$ws4.Range("F31").Value = 99
$ws4.Range("F33").Value = 6537
$ws4.Range("F34").Value = 95
$ws4.Range("F37").Value = 567
$ws4.Range("F38").Value = 1959
$ws4.Range("F41").Value = 89
$ws4.Range("F42").Value = 1019
$ws4.Range("F44").Value = 216
$ws4.Range("F46").Value = 123
$ws4.Range("F48").Value = 1742
